# Removing new user dependency from watch list test scripts
#
# The "Test Suite" sheet lists test rows keyed by a TSID (column A). Three
# of those TSIDs referred to a since-removed "new user" persona (F/C/D
# Suite). Re-point those rows at existing personas (Notifications /
# Authoring / Profile) instead, so the watch-list test scripts no longer
# depend on a user that doesn't exist anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Row 4: "C Suite" -> "Authoring" (TSID for the Authoring module test)
$ws.Range("A4").Value = "Authoring"

# Row 5: "D Suite" -> "Profile" (TSID for the Profile module test)
$ws.Range("A5").Value = "Profile"

# Row 7: "F Suite" -> "Notifications" (TSID for the Notification module test)
$ws.Range("A7").Value = "Notifications"

# Move the active selection to reflect where the editor ended up (A7).
$ws.Activate() | Out-Null
$ws.Range("A7").Select() | Out-Null
